$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 658
$ws.Range("I2").Value = 1797
$ws.Range("J2").Value = 7529
$ws.Range("K2").Value = 35
$ws.Range("L2").Value = 2136
$ws.Range("N2").Value = 1210
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 103
$ws.Range("S2").Value = 851
$ws.Range("T2").Value = 1265
$ws.Range("U2").Value = 95
$ws.Range("V2").Value = 11546
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 11671
$ws.Range("Y2").Value = 16
$ws.Range("Z2").Value = 160
$ws.Range("AA2").Value = 63
